$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rebuild column A (row index) as a sequential 1..200 counter.
$ws.Range("A2").Value = 1
$ws.Range("A3").Formula = "=A2+1"
$ws.Range("A4:A67").Formula = "=A3+1"
$ws.Range("A68:A131").Formula = "=A67+1"
$ws.Range("A132:A195").Formula = "=A131+1"
$ws.Range("A196:A201").Formula = "=A195+1"

# Update selection to match the saved view state.
$ws.Activate()
$ws.Range("A3:A201").Select()
